$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '54.015.77'
$ws.Range('E2').Value = '  -5.92%  '

# Row 3
$ws.Range('D3').Value = '2.841.77'
$ws.Range('E3').Value = '  -9.86%  '

# Row 4
$ws.Range('E4').Value = '  -0.15%  '

# Row 5
$ws.Range('D5').Value = '''459.71'
$ws.Range('E5').Value = '  -12.87%  '

# Row 6
$ws.Range('D6').Value = '''123.49'
$ws.Range('E6').Value = '  -8.32%  '

# Row 7
$ws.Range('E7').Value = '  -0.21%  '

# Row 8
$ws.Range('D8').Value = '2.858.22'
$ws.Range('E8').Value = '  -8.96%  '

# Row 9
$ws.Range('D9').Value = '''0.401'
$ws.Range('E9').Value = '  -11.03%  '

# Row 10
$ws.Range('D10').Value = '''6.56'
$ws.Range('E10').Value = '  -9.18%  '

# Row 11
$ws.Range('D11').Value = '''0.0952'
$ws.Range('E11').Value = '  -14.35%  '

# Row 12
$ws.Range('D12').Value = '''0.329'
$ws.Range('E12').Value = '  -16.33%  '

# Row 13
$ws.Range('D13').Value = '''0.122'
$ws.Range('E13').Value = '  -4.80%  '

# Row 14
$ws.Range('D14').Value = '3.326.01'
$ws.Range('E14').Value = '  -9.97%  '

# Row 15
$ws.Range('D15').Value = '''23.01'
$ws.Range('E15').Value = '  -11.06%  '

# Row 16
$ws.Range('D16').Value = '54.011.76'
$ws.Range('E16').Value = '  -6.14%  '

# Row 17
$ws.Range('D17').Value = '2.850.65'
$ws.Range('E17').Value = '  -9.52%  '

# Row 18
$ws.Range('D18').Value = '''0.0000132'
$ws.Range('E18').Value = '  -13.99%  '

# Row 19
$ws.Range('D19').Value = '''5.31'
$ws.Range('E19').Value = '  -9.24%  '

# Row 20
$ws.Range('D20').Value = '''11.29'
$ws.Range('E20').Value = '  -14.30%  '

# Row 21
$ws.Range('D21').Value = '''6.96'
$ws.Range('E21').Value = '  -13.30%  '

# Row 22
$ws.Range('D22').Value = '''294.69'
$ws.Range('E22').Value = '  -15.85%  '

# Row 23
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.23%  '

# Row 24
$ws.Range('D24').Value = '''0.436'
$ws.Range('E24').Value = '  -15.00%  '

# Row 25
$ws.Range('D25').Value = '''57.94'
$ws.Range('E25').Value = '  -16.76%  '

# Row 26
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.22%  '

# Row 27
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '''0.149'
$ws.Range('E27').Value = '  -10.34%  '

# Row 28
$ws.Range('B28').Value = 'USDe'
$ws.Range('C28').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D28').Value = '''0.998'
$ws.Range('E28').Value = '  -0.08%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0786'
$ws.Range('E29').Value = '  -18.09%  '

# Row 30
$ws.Range('D30').Value = '''6.03'
$ws.Range('E30').Value = '  -12.38%  '

# Row 31
$ws.Range('D31').Value = '''6.13'
$ws.Range('E31').Value = '  -11.96%  '

# Row 32
$ws.Range('D32').Value = '''1.08'
$ws.Range('E32').Value = '  -10.92%  '

# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''18.54'
$ws.Range('E33').Value = '  -14.06%  '

# Row 34
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '''1.57'
$ws.Range('E34').Value = '  -16.20%  '

# Row 35
$ws.Range('D35').Value = '''136.51'
$ws.Range('E35').Value = '  -13.68%  '

# Row 36
$ws.Range('D36').Value = '''4.11'
$ws.Range('E36').Value = '  -16.51%  '

# Row 37
$ws.Range('D37').Value = '''5.34'
$ws.Range('E37').Value = '  -14.30%  '

# Row 38
$ws.Range('D38').Value = '''1.20'
$ws.Range('E38').Value = '  -15.55%  '

# Row 39
$ws.Range('D39').Value = '''22.87'
$ws.Range('E39').Value = '  -12.06%  '

# Row 40
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '''0.0610'
$ws.Range('E40').Value = '  -12.89%  '

# Row 41
$ws.Range('E41').Value = '  -0.22%  '

# Row 42
$ws.Range('B42').Value = 'RenzoRestakedETH'
$ws.Range('C42').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D42').Value = '2.862.47'
$ws.Range('E42').Value = '  -9.94%  '

# Row 43
$ws.Range('D43').Value = '''34.83'
$ws.Range('E43').Value = '  -13.42%  '

# Row 44
$ws.Range('E44').Value = '  -13.88%  '

# Row 45
$ws.Range('D45').Value = '''0.922'
$ws.Range('E45').Value = '  -15.74%  '

# Row 46
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '''3.38'
$ws.Range('E46').Value = '  -14.68%  '

# Row 47
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '''1.29'
$ws.Range('E47').Value = '  -11.96%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.020.19'
$ws.Range('E48').Value = '  -10.60%  '

# Row 49
$ws.Range('D49').Value = '''5.30'
$ws.Range('E49').Value = '  -14.56%  '

# Row 50
$ws.Range('D50').Value = '''0.0212'
$ws.Range('E50').Value = '  -10.40%  '

# Row 51
$ws.Range('D51').Value = '''17.45'
$ws.Range('E51').Value = '  -15.15%  '
